$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of beverage-sale data appended below the existing single data row.
# Columns F (price) and G (date) look numeric/date-like but must be stored
# as plain text (matching the inlineStr cells used elsewhere in the sheet),
# so their number format is forced to "Text" before the value is assigned —
# this stops Excel's automatic type inference from turning "30.0" into the
# number 30 or "2024-09-20" into a date serial.
$ws.Range("F3").NumberFormat = "@"
$ws.Range("G3").NumberFormat = "@"

$ws.Range("A3").Value = "4b31a1f2-c211-481c-8fc7-8b89583eff35"
$ws.Range("B3").Value = "dayli"
$ws.Range("C3").Value = "Juices"
$ws.Range("D3").Value = "Coca-Cola"
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = "30.0"
$ws.Range("G3").Value = "2024-09-20"
$ws.Range("H3").Value = "18:22:22"
